# Auto-generated edit script: update crypto price/volume table (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'72.002.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.18%  "
$ws.Range("D3").Value = "'4.034.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.70%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'517.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").Value = "'148.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.52%  "
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "'0.733"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "'47.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +14.26%  "
$ws.Range("D13").Value = "'10.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.81%  "
$ws.Range("D14").Value = "'4.682.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.11%  "
$ws.Range("D15").Value = "'4.038.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.81%  "
$ws.Range("D16").Value = "'21.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.74%  "
$ws.Range("D17").Value = "'14.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("E19").Value = "  -2.30%  "
$ws.Range("D20").Value = "'72.040.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.54%  "
$ws.Range("D21").Value = "'435.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.00%  "
$ws.Range("D22").Value = "'97.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.73%  "
$ws.Range("D23").Value = "'3.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.98%  "
$ws.Range("D24").Value = "'14.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.36%  "
$ws.Range("D25").Value = "'11.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.19%  "
$ws.Range("D26").Value = "'4.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").Value = "'11.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.44%  "
$ws.Range("D28").Value = "'36.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.09%  "
$ws.Range("D29").Value = "'3.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.02%  "
$ws.Range("D30").Value = "'697.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").Value = "'13.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.18%  "
$ws.Range("D32").Value = "'0.128"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.03%  "
$ws.Range("D33").Value = "'7.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +19.06%  "
$ws.Range("D34").Value = "'68.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").Value = "'0.0₃0886"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.41%  "
$ws.Range("B36").Value = "TheGraph"
$ws.Range("C36").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D36").Value = "'0.436"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("B37").Value = "ThetaToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D37").Value = "'3.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +24.67%  "
$ws.Range("D38").Value = "'40.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("E39").Value = "  +3.81%  "
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "'0.0486"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.22%  "
$ws.Range("D43").Value = "'3.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.25%  "
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("E45").Value = "  +5.83%  "
$ws.Range("E46").Value = "  +3.71%  "
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("D48").Value = "'9.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.92%  "
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("D50").Value = "'0.000268"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +19.35%  "
$ws.Range("D51").Value = "'0.0₆0340"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.09%  "
